$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.840.09"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "'1.890.74"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'0.7815"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").Value = "'244.01"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").Value = "'0.9995"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.3141"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").Value = "'0.07323"
$ws.Range("E9").Value = "  +4.27%  "
$ws.Range("D10").Value = "'25.36"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("D11").Value = "'0.08134"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").Value = "'0.7662"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "'5.468"
$ws.Range("E13").Value = "  +2.97%  "
$ws.Range("D14").Value = "'1.881.85"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "'93.11"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "'6.207"
$ws.Range("E16").Value = "  +4.93%  "
$ws.Range("D17").Value = "'29.835.90"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "'13.94"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").Value = "'245.46"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").Value = "'0.000007878"
$ws.Range("E20").Value = "  +2.65%  "
$ws.Range("D21").Value = "'0.9994"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "'8.148"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").Value = "'2.130.50"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").Value = "'0.9995"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "'0.1590"
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("D26").Value = "'9.459"
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("D27").Value = "'161.82"
$ws.Range("E27").Value = "  -1.14%  "
$ws.Range("D28").Value = "'18.77"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").Value = "'2.034"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").Value = "'1.452"
$ws.Range("E30").Value = "  +5.55%  "
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("D32").Value = "'4.477"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("D33").Value = "'0.05594"
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("D34").Value = "'4.082"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "'1.253"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").Value = "'0.7552"
$ws.Range("E36").Value = "  +2.59%  "
$ws.Range("D37").Value = "'0.9960"
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("D38").Value = "'2.638"
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("D39").Value = "'0.01935"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("D40").Value = "'2.778"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").Value = "'1.143.61"
$ws.Range("E41").Value = "  +10.86%  "
$ws.Range("D42").Value = "'0.4450"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("D43").Value = "'73.84"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("D44").Value = "'5.964"
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("D45").Value = "'0.8560"
$ws.Range("E45").Value = "  +2.04%  "
$ws.Range("D46").Value = "'0.9995"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "'1.898"
$ws.Range("E47").Value = "  +2.37%  "
$ws.Range("D50").Value = "'9.807"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").Value = "'7.513"
$ws.Range("E51").Value = "  +0.95%  "

# Row 48/49 swap (Quant <-> SynthetixNetwork)
$ws.Range("B48").Value = "SynthetixNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D48").Value = "'3.118"
$ws.Range("E48").Value = "  +7.07%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'101.92"
$ws.Range("E49").Value = "  -0.12%  "
